# Updated cryptos list on Fri Jul 26 05:46:38 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the value to be stored as text even if it looks like a number
    # (e.g. "1.00", "6.72"), then strip the quote-prefix formatting that
    # Excel applies so the cell keeps the default (unstyled) look.
    $range.Value = "'" + $text
    $range.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "66.914.61"
$ws.Range("E2").Value = "  +4.07%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.256.84"
$ws.Range("E3").Value = "  +2.12%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "579.43"
$ws.Range("E5").Value = "  +2.83%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "176.66"
$ws.Range("E6").Value = "  +3.21%  "

# Row 7 - was XRP, now USDC
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - was USDC, now XRP
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D8") "0.605"
$ws.Range("E8").Value = "  +0.62%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws.Range("D9") "3.256.88"
$ws.Range("E9").Value = "  +2.20%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +4.25%  "

# Row 11 - Toncoin
Set-TextValue $ws.Range("D11") "6.72"
$ws.Range("E11").Value = "  +1.12%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("D12") "0.408"
$ws.Range("E12").Value = "  +2.95%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "3.820.31"
$ws.Range("E13").Value = "  +2.01%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.10%  "

# Row 15 - Avalanche
Set-TextValue $ws.Range("D15") "28.06"
$ws.Range("E15").Value = "  +1.68%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "66.911.07"
$ws.Range("E16").Value = "  +4.08%  "

# Row 17 - ShibaInu
Set-TextValue $ws.Range("D17") "0.0000167"
$ws.Range("E17").Value = "  +2.84%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "3.255.48"
$ws.Range("E18").Value = "  +2.13%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +2.24%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "13.40"
$ws.Range("E20").Value = "  +2.25%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "369.09"
$ws.Range("E21").Value = "  +4.64%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "7.60"
$ws.Range("E22").Value = "  +5.58%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.15%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "70.77"
$ws.Range("E24").Value = "  +1.85%  "

# Row 25 - Polygon
Set-TextValue $ws.Range("D25") "0.509"
$ws.Range("E25").Value = "  +0.73%  "

# Row 26 - was PEPE, now WrappedeETH
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D26") "3.392.99"
$ws.Range("E26").Value = "  +2.17%  "

# Row 27 - was WrappedeETH, now PEPE
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D27") "0.0000119"
$ws.Range("E27").Value = "  -0.01%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D28") "9.77"
$ws.Range("E28").Value = "  +1.66%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  +2.05%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  -0.10%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +4.33%  "

# Row 32 - NEARProtocol
Set-TextValue $ws.Range("D32") "5.65"
$ws.Range("E32").Value = "  -0.40%  "

# Row 33 - EthereumClassic
Set-TextValue $ws.Range("D33") "22.55"
$ws.Range("E33").Value = "  +1.86%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  -0.13%  "

# Row 35 - was Monero, now Fetch.AI
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D35") "1.24"
$ws.Range("E35").Value = "  +2.95%  "

# Row 36 - was Fetch.AI, now Aptos
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D36") "6.77"
$ws.Range("E36").Value = "  +1.71%  "

# Row 37 - was Aptos, now Monero
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D37") "170.78"
$ws.Range("E37").Value = "  +9.32%  "

# Row 38 - ImmutableX
Set-TextValue $ws.Range("D38") "1.51"
$ws.Range("E38").Value = "  +4.41%  "

# Row 39 - Mantle
Set-TextValue $ws.Range("D39") "0.858"
$ws.Range("E39").Value = "  +5.84%  "

# Row 40 - Stacks
Set-TextValue $ws.Range("D40") "1.84"
$ws.Range("E40").Value = "  +9.26%  "

# Row 41 - EnergySwap
Set-TextValue $ws.Range("D41") "27.07"
$ws.Range("E41").Value = "  +4.24%  "

# Row 42 - dogwifhat
Set-TextValue $ws.Range("D42") "2.57"
$ws.Range("E42").Value = "  +2.71%  "

# Row 43 - Maker
Set-TextValue $ws.Range("D43") "2.747.63"
$ws.Range("E43").Value = "  +3.16%  "

# Row 44 - RenderToken
Set-TextValue $ws.Range("D44") "6.42"
$ws.Range("E44").Value = "  +6.58%  "

# Row 45 - Filecoin
$ws.Range("E45").Value = "  +3.68%  "

# Row 46 - Bittensor
Set-TextValue $ws.Range("D46") "341.91"
$ws.Range("E46").Value = "  +3.87%  "

# Row 47 - OKB
Set-TextValue $ws.Range("D47") "40.34"
$ws.Range("E47").Value = "  +4.47%  "

# Row 48 - Hedera
Set-TextValue $ws.Range("D48") "0.0673"
$ws.Range("E48").Value = "  +3.11%  "

# Row 49 - InjectiveProtocol
Set-TextValue $ws.Range("D49") "24.72"
$ws.Range("E49").Value = "  +4.46%  "

# Row 50 - VeChain
Set-TextValue $ws.Range("D50") "0.0278"
$ws.Range("E50").Value = "  +2.48%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  +2.52%  "
